# Rename transcript speaker codes in the DataSheet.
# "RT1"/"RT5" (teacher/researcher codes) -> "T"
# "Class" (whole class speaking) -> "SS"
# Also fix up a stray reference to "class" inside a sentence cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows whose Speaker column (D) currently holds "RT1" or "RT5" and must become "T"
$teacherRows = @(
    2,3,4,5,6,7,8,11,13,14,15,16,17,18,19,20,21,22,23,24,26,29,30,31,
    33,34,35,36,37,38,39,42,50,51,52,53,54,55,56,57,58,59,64,65,67,
    70,71,72,74,76,77,78,79,83,85,86,87,88,89,90,91,92,93,94,95,96,
    97,98,99,100,101,102,103,104,105,107,108,109,110,111,112,113,114,
    115,116,117,118,119,120,121,122,123,125,126,127,139,140,141,147,
    148,149,151,154,155,159,160,164,165,166,167,168,169,170,172,174,
    175,180,181,184,185,187,188,189,190,191,192,193
)

foreach ($r in $teacherRows) {
    $ws.Cells.Item($r, 4).Value = "T"
}

# Rows whose Speaker column (D) currently holds "Class" and must become "SS"
$classRows = @(25, 27, 32)

foreach ($r in $classRows) {
    $ws.Cells.Item($r, 4).Value = "SS"
}

# Fix the inline reference to "class" inside row 189's sentence text
$ws.Cells.Item(189, 5).Value = "[to SS]  I'd like you to  think about the little number line you made, the fraction  number line between zero and one."
